$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# 1) TS sheet: add a new row (row 6) to Table2 with the "TS04" strategy
# ----------------------------------------------------------------------
$wsTS = $wb.Worksheets.Item("TS")
$loTS = $wsTS.ListObjects.Item("Table2")
$newTSRow = $loTS.ListRows.Add()

$wsTS.Range("A6").Value = "TS04"
$wsTS.Range("B6").Value = "future 202107, con pandemia"
$wsTS.Range("C6").Value = 202107
$wsTS.Range("D6").Value = "202105, 202104, 202103, 202102, 202101, 202012, 202011, 202010, 202009, 202008, 202007, 202006, 202005, 202004, 202003, 202002, 202001, 201912, 201911, 201910, 201909, 201908, 201907"
$wsTS.Range("E6").Value = "202103, 202102, 202101, 202012, 202011, 202010, 202009, 202008, 202007, 202006, 202005, 202004, 202003, 202002, 202001, 201912, 201911, 201910, 201909, 201908, 201907, 201906, 201905"
$wsTS.Range("F6").Value = 202104
$wsTS.Range("G6").Value = 202105
$wsTS.Range("H6").Value = 0.5
$wsTS.Range("I6").Value = 10881

$wsTS.Activate()
$wsTS.Range("C7").Select()

# ----------------------------------------------------------------------
# 2) FE sheet: rename current "FE02" row to "FE01.02" and add three new
#    rows (FE02, FE03, FE04) to Table1
# ----------------------------------------------------------------------
$wsFE = $wb.Worksheets.Item("FE")

$wsFE.Range("A3").Value = "FE01.02"

$loFE = $wsFE.ListObjects.Item("Table1")

$row4 = $loFE.ListRows.Add()
$wsFE.Range("A4").Value = "FE02"
$wsFE.Range("B4").Value = 200
$wsFE.Range("C4").Value = 7
$wsFE.Range("D4").Value = 600
$wsFE.Range("E4").Value = 50
$wsFE.Range("F4").Value = 10881
$wsFE.Range("G4").Value = 0.999
$wsFE.Range("H4").Value = $true
$wsFE.Range("I4").Value = $true
$wsFE.Range("J4").Value = $true
$wsFE.Range("K4").Value = "Todo TRUE"
$wsFE.Range("L4").Value = "Todo TRUE"

$row5 = $loFE.ListRows.Add()
$wsFE.Range("A5").Value = "FE03"
$wsFE.Range("B5").Value = 200
$wsFE.Range("C5").Value = 6
$wsFE.Range("D5").Value = 600
$wsFE.Range("E5").Value = 50
$wsFE.Range("F5").Value = 10881
$wsFE.Range("G5").Value = 0.999
$wsFE.Range("H5").Value = $true
$wsFE.Range("I5").Value = $true
$wsFE.Range("J5").Value = $true
$wsFE.Range("K5").Value = "Todo TRUE"
$wsFE.Range("L5").Value = "Todo TRUE"

$row6 = $loFE.ListRows.Add()
$wsFE.Range("A6").Value = "FE04"
$wsFE.Range("B6").Value = 300
$wsFE.Range("C6").Value = 5
$wsFE.Range("D6").Value = 600
$wsFE.Range("E6").Value = 50
$wsFE.Range("F6").Value = 10881
$wsFE.Range("G6").Value = 0.999
$wsFE.Range("H6").Value = $true
$wsFE.Range("I6").Value = $true
$wsFE.Range("J6").Value = $true
$wsFE.Range("K6").Value = "Todo TRUE"
$wsFE.Range("L6").Value = "Todo TRUE"

$wsFE.Activate()
$wsFE.Range("A7").Select()

$wb.Save()
